$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "47.468.04"; E = "  +4.80%  " },
    @{ Row = 3;  D = "2.483.85";  E = "  +2.48%  " },
    @{ Row = 4;  D = $null;       E = "  +0.07%  " },
    @{ Row = 5;  D = "322.78";    E = "  +1.44%  " },
    @{ Row = 6;  D = "104.92";    E = "  +2.44%  " },
    @{ Row = 7;  D = $null;       E = "  +1.89%  " },
    @{ Row = 8;  D = $null;       E = "  +0.02%  " },
    @{ Row = 9;  D = $null;       E = "  +2.28%  " },
    @{ Row = 10; D = "38.05";     E = "  +7.00%  " },
    @{ Row = 11; D = $null;       E = "  +1.21%  " },
    @{ Row = 12; D = $null;       E = "  +1.13%  " },
    @{ Row = 13; D = "18.21";     E = "  +0.22%  " },
    @{ Row = 14; D = "7.14";      E = "  +1.73%  " },
    @{ Row = 15; D = "2.874.21";  E = "  +2.48%  " },
    @{ Row = 16; D = "2.486.52";  E = "  +2.49%  " },
    @{ Row = 17; D = "0.841";     E = "  -0.15%  " },
    @{ Row = 18; D = "47.366.77"; E = "  +4.80%  " },
    @{ Row = 19; D = "12.63";     E = "  +3.48%  " },
    @{ Row = 20; D = $null;       E = "  +3.38%  " },
    @{ Row = 21; D = $null;       E = "  +1.70%  " },
    @{ Row = 22; D = "70.66";     E = "  +2.73%  " },
    @{ Row = 23; D = "250.67";    E = "  +2.72%  " },
    @{ Row = 24; D = "2.39";      E = "  +6.00%  " },
    @{ Row = 25; D = $null;       E = "  +2.78%  " },
    @{ Row = 26; D = "26.09";     E = "  +2.17%  " },
    @{ Row = 27; D = $null;       E = "  -0.07%  " },
    @{ Row = 28; D = "9.98";      E = "  +4.29%  " },
    @{ Row = 29; D = "2.25";      E = "  +8.68%  " },
    @{ Row = 30; D = "34.92";     E = "  +6.51%  " },
    @{ Row = 31; D = "0.134";     E = "  +6.51%  " },
    @{ Row = 32; D = "49.42";     E = "  +0.67%  " },
    @{ Row = 33; D = "19.82";     E = "  -1.86%  " },
    @{ Row = 34; D = "5.35";      E = "  +2.81%  " },
    @{ Row = 35; D = $null;       E = "  +2.10%  " },
    @{ Row = 36; D = $null;       E = "  +0.18%  " },
    @{ Row = 37; D = $null;       E = "  +3.65%  " },
    @{ Row = 38; D = $null;       E = "  +3.24%  " },
    @{ Row = 39; D = $null;       E = "  +4.22%  " },
    @{ Row = 40; D = $null;       E = "  +1.87%  " },
    @{ Row = 41; D = $null;       E = "  +1.75%  " },
    @{ Row = 42; D = "121.24";    E = "  -4.28%  " },
    @{ Row = 43; D = "21.11";     E = "  +3.37%  " },
    @{ Row = 44; D = $null;       E = "  +2.57%  " },
    @{ Row = 45; D = "1.962.34";  E = "  +1.56%  " },
    @{ Row = 46; D = "2.97";      E = "  +1.63%  " },
    @{ Row = 48; D = $null;       E = "  +1.05%  " },
    @{ Row = 49; D = $null;       E = "  -1.48%  " },
    @{ Row = 50; D = "5.29";      E = "  +12.22%  " },
    @{ Row = 51; D = "79.36";     E = "  +3.77%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    $eCell = $ws.Cells.Item($u.Row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $u.E
}
